$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 0.1857707509881423
$ws.Range("C2").Value = 0.549407114624506
$ws.Range("J2").Value = 0.01383399209486166
$ws.Range("O2").Value = 0.001976284584980237
$ws.Range("P2").Value = 0.1482213438735178
$ws.Range("S2").Value = 0.1007905138339921
# Row 3
$ws.Range("B3").Value = 0.02033898305084746
$ws.Range("C3").Value = 0.04067796610169491
$ws.Range("J3").Value = 0.02033898305084746
$ws.Range("P3").Value = 0.7762711864406779
$ws.Range("S3").Value = 0.1423728813559322
# Row 4
$ws.Range("J4").Value = 0.04761904761904762
$ws.Range("P4").Value = 0.6547619047619048
$ws.Range("S4").Value = 0.2976190476190476
# Row 5
$ws.Range("P5").Value = 0.3333333333333333
$ws.Range("S5").Value = 0.6666666666666666
# Row 6
$ws.Range("B6").Value = 0.06666666666666667
$ws.Range("D6").Value = 0.002469135802469136
$ws.Range("E6").Value = 0.004938271604938272
$ws.Range("F6").Value = 0.04197530864197531
$ws.Range("J6").Value = 0.254320987654321
$ws.Range("O6").Value = 0.02469135802469136
$ws.Range("Q6").Value = 0.1654320987654321
$ws.Range("R6").Value = 0.07407407407407407
$ws.Range("S6").Value = 0.3654320987654321
# Row 7
$ws.Range("B7").Value = 0.08900523560209424
$ws.Range("D7").Value = 0.03403141361256545
$ws.Range("F7").Value = 0.06806282722513089
$ws.Range("J7").Value = 0.1465968586387434
$ws.Range("O7").Value = 0.01570680628272251
$ws.Range("Q7").Value = 0.1675392670157068
$ws.Range("R7").Value = 0.06282722513089005
$ws.Range("S7").Value = 0.4162303664921466
# Row 8
$ws.Range("B8").Value = 0.1082887700534759
$ws.Range("D8").Value = 0.01871657754010695
$ws.Range("F8").Value = 0.06283422459893048
$ws.Range("J8").Value = 0.08689839572192513
$ws.Range("O8").Value = 0.01871657754010695
$ws.Range("Q8").Value = 0.213903743315508
$ws.Range("R8").Value = 0.08422459893048129
$ws.Range("S8").Value = 0.4064171122994653
# Row 9
$ws.Range("B9").Value = 0.08732394366197183
$ws.Range("D9").Value = 0.01690140845070422
$ws.Range("F9").Value = 0.07887323943661972
$ws.Range("J9").Value = 0.1042253521126761
$ws.Range("O9").Value = 0.03098591549295775
$ws.Range("Q9").Value = 0.1633802816901408
$ws.Range("R9").Value = 0.08450704225352113
$ws.Range("S9").Value = 0.4338028169014084
# Row 10
$ws.Range("B10").Value = 0.1042039355992844
$ws.Range("D10").Value = 0.02549194991055456
$ws.Range("E10").Value = 0.0004472271914132379
$ws.Range("F10").Value = 0.06618962432915922
$ws.Range("J10").Value = 0.1122540250447227
$ws.Range("O10").Value = 0.01967799642218247
$ws.Range("Q10").Value = 0.2101967799642218
$ws.Range("R10").Value = 0.07379248658318426
$ws.Range("S10").Value = 0.3877459749552773
# Row 11
$ws.Range("G11").Value = 0.1468646864686469
$ws.Range("J11").Value = 0.08250825082508251
$ws.Range("K11").Value = 0.202970297029703
$ws.Range("L11").Value = 0.5561056105610561
$ws.Range("S11").Value = 0.01155115511551155
# Row 12
$ws.Range("G12").Value = 0.6887052341597796
$ws.Range("J12").Value = 0.2286501377410468
$ws.Range("K12").Value = 0.01101928374655647
$ws.Range("L12").Value = 0.02479338842975207
$ws.Range("S12").Value = 0.04683195592286502
# Row 13
$ws.Range("G13").Value = 0.5730337078651685
$ws.Range("J13").Value = 0.3370786516853932
$ws.Range("S13").Value = 0.0898876404494382
# Row 14
$ws.Range("G14").Value = 0.6
$ws.Range("J14").Value = 0.2
$ws.Range("S14").Value = 0.2
# Row 15
$ws.Range("F15").Value = 0.03073286052009456
$ws.Range("H15").Value = 0.1205673758865248
$ws.Range("I15").Value = 0.08037825059101655
$ws.Range("J15").Value = 0.3782505910165485
$ws.Range("K15").Value = 0.05437352245862884
$ws.Range("M15").Value = 0.009456264775413711
$ws.Range("O15").Value = 0.06146572104018912
$ws.Range("S15").Value = 0.2647754137115839
# Row 16
$ws.Range("F16").Value = 0.01994301994301994
$ws.Range("H16").Value = 0.150997150997151
$ws.Range("I16").Value = 0.09401709401709402
$ws.Range("J16").Value = 0.3846153846153846
$ws.Range("K16").Value = 0.1111111111111111
$ws.Range("M16").Value = 0.005698005698005698
$ws.Range("O16").Value = 0.05698005698005698
$ws.Range("S16").Value = 0.1766381766381766
# Row 17
$ws.Range("F17").Value = 0.02085889570552147
$ws.Range("H17").Value = 0.1779141104294479
$ws.Range("I17").Value = 0.09570552147239264
$ws.Range("J17").Value = 0.3938650306748466
$ws.Range("K17").Value = 0.1092024539877301
$ws.Range("M17").Value = 0.0245398773006135
$ws.Range("N17").Value = 0.001226993865030675
$ws.Range("O17").Value = 0.06257668711656442
$ws.Range("S17").Value = 0.1141104294478528
# Row 18
$ws.Range("F18").Value = 0.02912621359223301
$ws.Range("H18").Value = 0.1812297734627832
$ws.Range("I18").Value = 0.07766990291262135
$ws.Range("J18").Value = 0.3818770226537217
$ws.Range("K18").Value = 0.1035598705501618
$ws.Range("M18").Value = 0.0226537216828479
$ws.Range("O18").Value = 0.06796116504854369
$ws.Range("S18").Value = 0.1359223300970874
# Row 19
$ws.Range("F19").Value = 0.01941747572815534
$ws.Range("H19").Value = 0.1899535669058675
$ws.Range("I19").Value = 0.07978049810046434
$ws.Range("J19").Value = 0.3541578725200507
$ws.Range("K19").Value = 0.122836639932461
$ws.Range("M19").Value = 0.02448290417897847
$ws.Range("N19").Value = 0.00295483326298016
$ws.Range("O19").Value = 0.07429295061207261
$ws.Range("S19").Value = 0.13212325875897
